$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.939.84"
$ws.Range("E2").Value = '  -2.29%  '

$ws.Range("D3").Value = "'1.796.27"
$ws.Range("E3").Value = '  -0.35%  '

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'317.12"
$ws.Range("E5").Value = '  +0.05%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = "'0.5313"
$ws.Range("E7").Value = '  -2.81%  '

$ws.Range("D8").Value = "'0.3876"
$ws.Range("E8").Value = '  +2.89%  '

$ws.Range("D9").Value = "'0.07453"
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").Value = "'41.43"
$ws.Range("E10").Value = '  -2.26%  '

$ws.Range("E11").Value = '  -2.25%  '

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = '  -0.07%  '

$ws.Range("D13").Value = "'6.176"
$ws.Range("E13").Value = '  +0.48%  '

$ws.Range("D14").Value = "'7.441"
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("E15").Value = '  -1.19%  '

$ws.Range("D16").Value = "'1.796.62"
$ws.Range("E16").Value = '  -0.20%  '

$ws.Range("D17").Value = "'88.43"
$ws.Range("E17").Value = '  -1.96%  '

$ws.Range("D18").Value = "'0.00001060"
$ws.Range("E18").Value = '  -0.45%  '

$ws.Range("D19").Value = "'0.06559"
$ws.Range("E19").Value = '  +1.64%  '

$ws.Range("D20").Value = "'0.9998"
$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").Value = "'5.956"
$ws.Range("E22").Value = '  +0.63%  '

$ws.Range("D23").Value = "'27.959.54"
$ws.Range("E23").Value = '  -2.32%  '

$ws.Range("D24").Value = "'11.11"
$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").Value = "'2.091"
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").Value = "'157.23"
$ws.Range("E26").Value = '  -0.93%  '

$ws.Range("D27").Value = "'20.15"
$ws.Range("E27").Value = '  -1.38%  '

$ws.Range("D28").Value = "'1.999.05"
$ws.Range("E28").Value = '  -0.45%  '

$ws.Range("D29").Value = "'2.299"
$ws.Range("E29").Value = '  -2.12%  '

$ws.Range("D30").Value = "'121.98"
$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'0.1088"
$ws.Range("E31").Value = '  +2.57%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'1.100"
$ws.Range("E32").Value = '  -0.38%  '

$ws.Range("D33").Value = "'3.665"
$ws.Range("E33").Value = '  -0.46%  '

$ws.Range("D34").Value = "'5.503"
$ws.Range("E34").Value = '  -2.37%  '

$ws.Range("D35").Value = "'0.07003"
$ws.Range("E35").Value = '  +7.92%  '

$ws.Range("D36").Value = "'0.2204"
$ws.Range("E36").Value = '  -1.87%  '

$ws.Range("D37").Value = "'0.02275"
$ws.Range("E37").Value = '  -1.11%  '

$ws.Range("D38").Value = "'5.080"
$ws.Range("E38").Value = '  +1.16%  '

$ws.Range("D39").Value = "'8.383"
$ws.Range("E39").Value = '  -3.99%  '

$ws.Range("D40").Value = "'11.23"
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("D41").Value = "'1.188"
$ws.Range("E41").Value = '  -1.49%  '

$ws.Range("D42").Value = "'0.6109"
$ws.Range("E42").Value = '  -1.84%  '

$ws.Range("D43").Value = "'1.419"
$ws.Range("E43").Value = '  -0.98%  '

$ws.Range("D44").Value = "'13.28"
$ws.Range("E44").Value = '  +0.49%  '

$ws.Range("D45").Value = "'3.677"
$ws.Range("E45").Value = '  -0.38%  '

$ws.Range("D46").Value = "'0.5708"
$ws.Range("E46").Value = '  -2.25%  '

$ws.Range("D47").Value = "'124.81"
$ws.Range("E47").Value = '  -1.35%  '

$ws.Range("D48").Value = "'1.178"
$ws.Range("E48").Value = '  +2.03%  '

$ws.Range("E49").Value = '  -1.20%  '

$ws.Range("E50").Value = '  -1.18%  '

$ws.Range("D51").Value = "'0.00000000295"
$ws.Range("E51").Value = '  +27.50%  '
